# Table A.8.B monthly refresh: October 2016 YTD -> November 2016 YTD
# (EIA electricity monthly table_a_8_b.xlsx, 2017-01-31 update, chunk 7)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update the subtitle text (row 2) to reference the new reporting month ---
$ws.Range("A2").Value = "by End-Use Sector, Census Division, and State, Year-to-Date through November 2016"

# --- Refreshed relative-standard-error figures (row -> column -> new value) ---

# Row 4  - New England
$ws.Range("C4").Value = 0.49

# Row 5  - Connecticut
$ws.Range("B5").Value = 0.13

# Row 6  - Maine
$ws.Range("B6").Value = 0.19

# Row 7  - Massachusetts
$ws.Range("B7").Value = 0.27
$ws.Range("D7").Value = 7

# Row 8  - New Hampshire
$ws.Range("B8").Value = 0.2

# Row 12 - New Jersey
$ws.Range("F12").Value = 0.26

# Row 14 - Pennsylvania
$ws.Range("B14").Value = 0.13
$ws.Range("C14").Value = 0.39

# Row 15 - East North Central
$ws.Range("F15").Value = 0.4

# Row 16 - Illinois
$ws.Range("B16").Value = 0.25
$ws.Range("F16").Value = 1

# Row 18 - Michigan
$ws.Range("B18").Value = 0.16

# Row 20 - Wisconsin
$ws.Range("B20").Value = 0.32

# Row 21 - West North Central
$ws.Range("B21").Value = 0.27

# Row 24 - Minnesota
$ws.Range("B24").Value = 0.47

# Row 25 - Missouri
$ws.Range("B25").Value = 0.48

# Row 29 - South Atlantic
$ws.Range("B29").Value = 0.29
$ws.Range("C29").Value = 0.38

# Row 31 - District of Columbia
$ws.Range("C31").Value = 0.4
$ws.Range("F31").Value = 0.3

# Row 32 - Florida
$ws.Range("B32").Value = 0.37

# Row 34 - Maryland
$ws.Range("B34").Value = 0.15
$ws.Range("C34").Value = 0.45

# Row 37 - Virginia
$ws.Range("B37").Value = 0.46
$ws.Range("F37").Value = 1

# Row 43 - Tennessee
$ws.Range("B43").Value = 0.46

# Row 44 - West South Central
$ws.Range("B44").Value = 0.39

# Row 48 - Texas
$ws.Range("B48").Value = 0.42

# Row 49 - Mountain
$ws.Range("B49").Value = 0.16

# Row 50 - Arizona
$ws.Range("B50").Value = 0.17
$ws.Range("D50").Value = 4

# Row 51 - Colorado
$ws.Range("B51").Value = 1

# Row 52 - Idaho
$ws.Range("B52").Value = 0.41

# Row 54 - Nevada
$ws.Range("B54").Value = 0.17

# Row 56 - Utah
$ws.Range("B56").Value = 1

# Row 59 - California
$ws.Range("B59").Value = 0.13

# Row 60 - Oregon
$ws.Range("B60").Value = 0.39

# Row 61 - Washington
$ws.Range("B61").Value = 0.33
